$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 - Net Receivables (Balance Sheet): G,H,I,J switch from 0 to "NA"
$ws.Range("G43").Value = "NA"
$ws.Range("H43").Value = "NA"
$ws.Range("I43").Value = "NA"
$ws.Range("J43").Value = "NA"

# Row 47 - Long Term Investments (Balance Sheet): D 1200->0, E:J "NA"->0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0

# Row 52 - Other Assets (Balance Sheet): D "NA"->1200
$ws.Range("D52").Value = 1200

# Row 57 - Accounts Payable (Balance Sheet): D 0->200, G:J 0->"NA"
$ws.Range("D57").Value = 200
$ws.Range("G57").Value = "NA"
$ws.Range("H57").Value = "NA"
$ws.Range("I57").Value = "NA"
$ws.Range("J57").Value = "NA"

# Row 58 - Short/Current Long Term Debt (Balance Sheet): D 0->1100, E 0->300, G:J 0->"NA"
$ws.Range("D58").Value = 1100
$ws.Range("E58").Value = 300
$ws.Range("G58").Value = "NA"
$ws.Range("H58").Value = "NA"
$ws.Range("I58").Value = "NA"
$ws.Range("J58").Value = "NA"

# Row 59 - Other Current Liabilities (Balance Sheet): D 1500->200, E 400->100
$ws.Range("D59").Value = 200
$ws.Range("E59").Value = 100

# Row 91 - Capital Expenditures (Cash Flow Statement): E -100->0
$ws.Range("E91").Value = 0
